$d = $word.ActiveDocument

$d.Content.Find.Execute("91÷8=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=7, 4", 2) | Out-Null
$d.Content.Find.Execute("60÷8=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "10÷3=3, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷5=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "86÷9=9, 5", 2) | Out-Null
$d.Content.Find.Execute("18÷2=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "83÷2=41, 1", 2) | Out-Null
$d.Content.Find.Execute("87÷5=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "27÷5=5, 2", 2) | Out-Null
$d.Content.Find.Execute("72÷2=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷7=11, 1", 2) | Out-Null
$d.Content.Find.Execute("55÷7=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "43÷9=4, 7", 2) | Out-Null
$d.Content.Find.Execute("59÷5=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=3, 2", 2) | Out-Null
$d.Content.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=13, 0", 2) | Out-Null
$d.Content.Find.Execute("10÷4=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "19÷7=2, 5", 2) | Out-Null
$d.Content.Find.Execute("63÷5=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "26÷2=13, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷5=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "86÷5=17, 1", 2) | Out-Null
$d.Content.Find.Execute("28÷2=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷9=8, 7", 2) | Out-Null
$d.Content.Find.Execute("88÷9=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "63÷6=10, 3", 2) | Out-Null
$d.Content.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷4=23, 2", 2) | Out-Null
$d.Content.Find.Execute("21÷4=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷9=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=25, 0", 2) | Out-Null
$d.Content.Find.Execute("53÷3=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "26÷4=6, 2", 2) | Out-Null
$d.Content.Find.Execute("15÷7=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=11, 0", 2) | Out-Null
$d.Content.Find.Execute("79÷4=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "63÷3=21, 0", 2) | Out-Null
$d.Content.Find.Execute("45÷8=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=49, 1", 2) | Out-Null
$d.Content.Find.Execute("85÷4=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "90÷6=15, 0", 2) | Out-Null
$d.Content.Find.Execute("22÷9=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "83÷3=27, 2", 2) | Out-Null
$d.Content.Find.Execute("34÷6=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "64÷7=9, 1", 2) | Out-Null
$d.Content.Find.Execute("35÷3=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "36÷4=9, 0", 2) | Out-Null
